$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price text parses as a plain number need to be forced to
# Text format first, so Excel keeps the exact original string (e.g. trailing
# zeros like "96.60") instead of auto-converting it to a numeric value. The
# style is then reset back to Normal so the cell's style index is unchanged.
$textPriceCells = @("D5","D6","D7","D10","D12","D14","D15","D19","D22","D24","D25","D27","D29","D33","D34","D36","D37","D41","D42","D43","D45","D46","D50")
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "37.821.03"
$ws.Range("E2").Value = "  -0.52%  "

$ws.Range("D3").Value = "2.030.18"
$ws.Range("E3").Value = "  -1.27%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "227.52"
$ws.Range("E5").Value = "  -1.13%  "

$ws.Range("D6").Value = "0.613"
$ws.Range("E6").Value = "  -0.54%  "

$ws.Range("D7").Value = "59.51"
$ws.Range("E7").Value = "  +2.22%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").Value = "  -0.65%  "

$ws.Range("D10").Value = "0.0810"
$ws.Range("E10").Value = "  +0.28%  "

$ws.Range("E11").Value = "  +0.08%  "

$ws.Range("D12").Value = "14.62"
$ws.Range("E12").Value = "  -0.11%  "

$ws.Range("D13").Value = "2.333.10"
$ws.Range("E13").Value = "  -1.14%  "

$ws.Range("D14").Value = "21.17"
$ws.Range("E14").Value = "  +2.39%  "

$ws.Range("D15").Value = "0.763"
$ws.Range("E15").Value = "  +1.20%  "

$ws.Range("E16").Value = "  -1.80%  "

$ws.Range("D17").Value = "2.031.32"
$ws.Range("E17").Value = "  -1.31%  "

$ws.Range("D18").Value = "37.739.55"
$ws.Range("E18").Value = "  -0.52%  "

$ws.Range("D19").Value = "6.03"
$ws.Range("E19").Value = "  -1.82%  "

$ws.Range("E20").Value = "  +0.25%  "

$ws.Range("E21").Value = "  -0.94%  "

$ws.Range("D22").Value = "224.84"
$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("D24").Value = "2.45"
$ws.Range("E24").Value = "  -0.31%  "

$ws.Range("D25").Value = "2.21"
$ws.Range("E25").Value = "  -1.42%  "

$ws.Range("E26").Value = "  -0.20%  "

$ws.Range("D27").Value = "165.32"
$ws.Range("E27").Value = "  -0.49%  "

$ws.Range("E28").Value = "  -3.07%  "

$ws.Range("D29").Value = "18.91"
$ws.Range("E29").Value = "  -0.76%  "

$ws.Range("E30").Value = "  -4.97%  "

$ws.Range("E31").Value = "  +0.61%  "

$ws.Range("E32").Value = "  -2.32%  "

$ws.Range("D33").Value = "2.08"
$ws.Range("E33").Value = "  +3.53%  "

$ws.Range("D34").Value = "4.49"
$ws.Range("E34").Value = "  -2.00%  "

$ws.Range("E35").Value = "  -1.67%  "

$ws.Range("D36").Value = "6.36"
$ws.Range("E36").Value = "  +6.43%  "

$ws.Range("D37").Value = "2.25"
$ws.Range("E37").Value = "  -3.78%  "

$ws.Range("E38").Value = "  -2.43%  "

$ws.Range("E39").Value = "  +0.04%  "

$ws.Range("D40").Value = "1.526.42"
$ws.Range("E40").Value = "  +2.74%  "

$ws.Range("D41").Value = "0.0219"
$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("D42").Value = "96.60"
$ws.Range("E42").Value = "  -1.89%  "

$ws.Range("D43").Value = "16.81"
$ws.Range("E43").Value = "  +0.84%  "

$ws.Range("E44").Value = "  -0.66%  "

$ws.Range("D45").Value = "0.0917"
$ws.Range("E45").Value = "  -2.07%  "

$ws.Range("D46").Value = "4.12"
$ws.Range("E46").Value = "  -1.29%  "

$ws.Range("E47").Value = "  -1.52%  "

$ws.Range("E48").Value = "  -1.05%  "

$ws.Range("E49").Value = "  -0.12%  "

$ws.Range("D50").Value = "7.09"
$ws.Range("E50").Value = "  +0.52%  "

$ws.Range("D51").Value = "2.220.83"
$ws.Range("E51").Value = "  -1.21%  "

# Restore the Normal style on the text-forced price cells so their style
# index matches the original (unstyled) cells again.
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).Style = "Normal"
}
